$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.507.61"
$ws.Range("E2").Value = "'  -0.32%  "

# Row 3
$ws.Range("D3").Value = "'1.731.49"
$ws.Range("E3").Value = "'  -0.68%  "

# Row 4
$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "'  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'247.62"
$ws.Range("E5").Value = "'  +0.60%  "

# Row 6
$ws.Range("E6").Value = "'  -0.03%  "

# Row 7
$ws.Range("D7").Value = "'0.4879"
$ws.Range("E7").Value = "'  +1.18%  "

# Row 8
$ws.Range("D8").Value = "'0.2675"
$ws.Range("E8").Value = "'  -0.97%  "

# Row 9
$ws.Range("D9").Value = "'0.06229"
$ws.Range("E9").Value = "'  -0.56%  "

# Row 10
$ws.Range("D10").Value = "'1.737.86"
$ws.Range("E10").Value = "'  -0.33%  "

# Row 11
$ws.Range("D11").Value = "'0.07068"
$ws.Range("E11").Value = "'  -0.85%  "

# Row 12
$ws.Range("D12").Value = "'15.67"
$ws.Range("E12").Value = "'  -1.24%  "

# Row 13
$ws.Range("D13").Value = "'4.622"
$ws.Range("E13").Value = "'  +2.16%  "

# Row 14
$ws.Range("D14").Value = "'0.6122"
$ws.Range("E14").Value = "'  -2.13%  "

# Row 15
$ws.Range("D15").Value = "'77.43"
$ws.Range("E15").Value = "'  -0.16%  "

# Row 16
$ws.Range("D16").Value = "'0.9995"
$ws.Range("E16").Value = "'  -0.05%  "

# Row 17
$ws.Range("D17").Value = "'26.506.12"
$ws.Range("E17").Value = "'  -0.33%  "

# Row 18
$ws.Range("D18").Value = "'0.9992"
$ws.Range("E18").Value = "'  -0.10%  "

# Row 19
$ws.Range("D19").Value = "'0.000007198"
$ws.Range("E19").Value = "'  +4.24%  "

# Row 20
$ws.Range("D20").Value = "'11.57"
$ws.Range("E20").Value = "'  -1.42%  "

# Row 21
$ws.Range("D21").Value = "'1.953.64"
$ws.Range("E21").Value = "'  -0.75%  "

# Row 22
$ws.Range("D22").Value = "'4.521"
$ws.Range("E22").Value = "'  -2.41%  "

# Row 23
$ws.Range("D23").Value = "'8.797"
$ws.Range("E23").Value = "'  -0.60%  "

# Row 24
$ws.Range("D24").Value = "'5.271"
$ws.Range("E24").Value = "'  -1.96%  "

# Row 25
$ws.Range("D25").Value = "'137.71"
$ws.Range("E25").Value = "'  +1.08%  "

# Row 26
$ws.Range("D26").Value = "'15.46"
$ws.Range("E26").Value = "'  +0.35%  "

# Row 27
$ws.Range("D27").Value = "'1.785"
$ws.Range("E27").Value = "'  -1.91%  "

# Row 28
$ws.Range("D28").Value = "'108.29"
$ws.Range("E28").Value = "'  +1.19%  "

# Row 29
$ws.Range("E29").Value = "'  -1.94%  "

# Row 30
$ws.Range("D30").Value = "'3.985"
$ws.Range("E30").Value = "'  -0.56%  "

# Row 31
$ws.Range("D31").Value = "'0.08019"
$ws.Range("E31").Value = "'  +1.62%  "

# Row 32
$ws.Range("D32").Value = "'3.696"
$ws.Range("E32").Value = "'  -1.47%  "

# Row 33
$ws.Range("D33").Value = "'0.04583"
$ws.Range("E33").Value = "'  -0.82%  "

# Row 34
$ws.Range("B34").Value = "'Frax"
$ws.Range("C34").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").Value = "'0.9991"
$ws.Range("E34").Value = "'  -0.06%  "

# Row 35
$ws.Range("B35").Value = "'HuobiToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.612"
$ws.Range("E35").Value = "'  -0.26%  "

# Row 36
$ws.Range("B36").Value = "'ARBITRUM"
$ws.Range("C36").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.007"
$ws.Range("E36").Value = "'  +0.77%  "

# Row 37
$ws.Range("B37").Value = "'ImmutableX"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.6363"
$ws.Range("E37").Value = "'  -0.86%  "

# Row 38
$ws.Range("B38").Value = "'TrustWalletToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'0.8981"
$ws.Range("E38").Value = "'  -4.66%  "

# Row 39
$ws.Range("B39").Value = "'RenderToken"
$ws.Range("C39").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.030"
$ws.Range("E39").Value = "'  +1.27%  "

# Row 40
$ws.Range("B40").Value = "'MXToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.390"
$ws.Range("E40").Value = "'  -1.53%  "

# Row 41
$ws.Range("B41").Value = "'PaxDollar"
$ws.Range("C41").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.004"
$ws.Range("E41").Value = "'  +0.11%  "

# Row 42
$ws.Range("B42").Value = "'VeChain"
$ws.Range("C42").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.01505"
$ws.Range("E42").Value = "'  -0.42%  "

# Row 43
$ws.Range("B43").Value = "'Quant"
$ws.Range("C43").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'101.60"
$ws.Range("E43").Value = "'  -10.39%  "

# Row 44
$ws.Range("B44").Value = "'FraxShare"
$ws.Range("C44").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.479"
$ws.Range("E44").Value = "'  -5.11%  "

# Row 45
$ws.Range("B45").Value = "'TheSandbox"
$ws.Range("C45").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.3907"
$ws.Range("E45").Value = "'  -0.40%  "

# Row 46
$ws.Range("B46").Value = "'Aptos"
$ws.Range("C46").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.019"
$ws.Range("E46").Value = "'  +4.18%  "

# Row 47
$ws.Range("B47").Value = "'Algorand"
$ws.Range("C47").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1185"
$ws.Range("E47").Value = "'  -3.26%  "

# Row 48
$ws.Range("B48").Value = "'Cronos"
$ws.Range("C48").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.05382"
$ws.Range("E48").Value = "'  +0.88%  "

# Row 49
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.942"
$ws.Range("E49").Value = "'  +0.26%  "

# Row 50
$ws.Range("B50").Value = "'Elrond"
$ws.Range("C50").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'30.64"
$ws.Range("E50").Value = "'  -0.47%  "

# Row 51
$ws.Range("B51").Value = "'NEARProtocol"
$ws.Range("C51").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.253"
$ws.Range("E51").Value = "'  -0.87%  "

